$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 2.66
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 3.25
$ws.Range("L2").Value = 1.5
$ws.Range("N2").Value = 3.3
$ws.Range("O2").Value = 1.43
$ws.Range("P2").Value = 1.76
$ws.Range("Q2").Value = 2.24
$ws.Range("R2").Value = 1.28
$ws.Range("S2").Value = 4.3
$ws.Range("T2").Value = 1.91
$ws.Range("V2").Value = 1.42
$ws.Range("W2").Value = 1.54
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 11.5
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 55
$ws.Range("AB2").Value = 9.6
$ws.Range("AD2").Value = 13.5
$ws.Range("AE2").Value = 40
$ws.Range("AH2").Value = 18.5
$ws.Range("AI2").Value = 90
$ws.Range("AJ2").Value = 38
$ws.Range("AK2").Value = 32
$ws.Range("AL2").Value = 48
$ws.Range("AN2").Value = 29
$ws.Range("AO2").Value = 42
$ws.Range("F3").Value = 1.62
$ws.Range("G3").Value = 1.67
$ws.Range("H3").Value = 6.6
$ws.Range("I3").Value = 7.6
$ws.Range("J3").Value = 3.85
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 1.56
$ws.Range("N3").Value = 2.98
$ws.Range("P3").Value = 1.65
$ws.Range("Q3").Value = 2.42
$ws.Range("T3").Value = 2.34
$ws.Range("U3").Value = 1.66
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 2.22
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 60
$ws.Range("AA3").Value = 290
$ws.Range("AB3").Value = 6.2
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 30
$ws.Range("AE3").Value = 150
$ws.Range("AF3").Value = 8.199999999999999
$ws.Range("AH3").Value = 32
$ws.Range("AI3").Value = 170
$ws.Range("AJ3").Value = 16
$ws.Range("AM3").Value = 580
$ws.Range("AN3").Value = 14.5
$ws.Range("AO3").Value = 1000
$ws.Range("F4").Value = 2.9
$ws.Range("G4").Value = 2.92
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 2.96
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 3
$ws.Range("Q4").Value = 2.42
$ws.Range("R4").Value = 1.25
$ws.Range("T4").Value = 1.91
$ws.Range("U4").Value = 1.96
$ws.Range("V4").Value = 1.51
$ws.Range("W4").Value = 1.52
$ws.Range("X4").Value = 10.5
$ws.Range("Z4").Value = 18
$ws.Range("AA4").Value = 48
$ws.Range("AC4").Value = 7
$ws.Range("AD4").Value = 15
$ws.Range("AF4").Value = 17.5
$ws.Range("AG4").Value = 14.5
$ws.Range("AJ4").Value = 48
$ws.Range("AK4").Value = 65
$ws.Range("AO4").Value = 42
$ws.Range("F5").Value = 1.63
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 6.6
$ws.Range("I5").Value = 7.4
$ws.Range("J5").Value = 3.95
$ws.Range("K5").Value = 4.1
$ws.Range("N5").Value = 3.2
$ws.Range("P5").Value = 1.76
$ws.Range("Q5").Value = 2.22
$ws.Range("R5").Value = 1.27
$ws.Range("T5").Value = 2.22
$ws.Range("U5").Value = 1.79
$ws.Range("V5").Value = 1.16
$ws.Range("W5").Value = 2.5
$ws.Range("X5").Value = 12
$ws.Range("Y5").Value = 18.5
$ws.Range("Z5").Value = 55
$ws.Range("AA5").Value = 230
$ws.Range("AB5").Value = 6.4
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 27
$ws.Range("AE5").Value = 130
$ws.Range("AF5").Value = 8.4
$ws.Range("AG5").Value = 10.5
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 140
$ws.Range("AJ5").Value = 16
$ws.Range("AK5").Value = 20
$ws.Range("AL5").Value = 85
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 200
$ws.Range("F6").Value = 1.72
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 5.8
$ws.Range("I6").Value = 6.4
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 3.9
$ws.Range("L6").Value = 1.47
$ws.Range("N6").Value = 3.35
$ws.Range("O6").Value = 1.39
$ws.Range("P6").Value = 1.8
$ws.Range("Q6").Value = 2.14
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 2.04
$ws.Range("U6").Value = 1.79
$ws.Range("V6").Value = 1.18
$ws.Range("W6").Value = 2.32
$ws.Range("X6").Value = 12.5
$ws.Range("Y6").Value = 19.5
$ws.Range("AA6").Value = 190
$ws.Range("AB6").Value = 7.6
$ws.Range("AC6").Value = 8.800000000000001
$ws.Range("AD6").Value = 25
$ws.Range("AF6").Value = 13
$ws.Range("AH6").Value = 24
$ws.Range("AJ6").Value = 17.5
$ws.Range("AK6").Value = 20
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 13.5
$ws.Range("F7").Value = 1.07
$ws.Range("G7").Value = 1.09
$ws.Range("J7").Value = 14.5
$ws.Range("K7").Value = 23
$ws.Range("N7").Value = 7.8
$ws.Range("P7").Value = 3.4
$ws.Range("R7").Value = 2.02
$ws.Range("S7").Value = 1.92
$ws.Range("T7").Value = 3.25
$ws.Range("U7").Value = 1.35
$ws.Range("AB7").Value = 1000
$ws.Range("AF7").Value = 9.6
$ws.Range("F8").Value = 1.67
$ws.Range("H8").Value = 6.4
$ws.Range("I8").Value = 7
$ws.Range("K8").Value = 3.95
$ws.Range("N8").Value = 2.98
$ws.Range("O8").Value = 1.48
$ws.Range("P8").Value = 1.67
$ws.Range("Q8").Value = 2.38
$ws.Range("R8").Value = 1.24
$ws.Range("T8").Value = 2.28
$ws.Range("U8").Value = 1.68
$ws.Range("W8").Value = 2.42
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 16.5
$ws.Range("AA8").Value = 250
$ws.Range("AB8").Value = 6.4
$ws.Range("AD8").Value = 28
$ws.Range("AF8").Value = 8.4
$ws.Range("AH8").Value = 30
$ws.Range("AI8").Value = 150
$ws.Range("AJ8").Value = 16.5
$ws.Range("AN8").Value = 16
$ws.Range("F9").Value = 1.51
$ws.Range("G9").Value = 1.56
$ws.Range("H9").Value = 7.4
$ws.Range("I9").Value = 8.800000000000001
$ws.Range("K9").Value = 4.6
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 3.55
$ws.Range("O9").Value = 1.36
$ws.Range("Q9").Value = 2.08
$ws.Range("R9").Value = 1.31
$ws.Range("S9").Value = 3.8
$ws.Range("T9").Value = 2.18
$ws.Range("U9").Value = 1.72
$ws.Range("V9").Value = 1.12
$ws.Range("W9").Value = 2.78
$ws.Range("X9").Value = 16
$ws.Range("AB9").Value = 6.8
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 38
$ws.Range("AE9").Value = 200
$ws.Range("AF9").Value = 8.6
$ws.Range("AM9").Value = 230
$ws.Range("AN9").Value = 10.5
$ws.Range("F10").Value = 1.42
$ws.Range("G10").Value = 1.43
$ws.Range("H10").Value = 10.5
$ws.Range("I10").Value = 13
$ws.Range("J10").Value = 4.5
$ws.Range("K10").Value = 5
$ws.Range("L10").Value = 1.38
$ws.Range("N10").Value = 4
$ws.Range("O10").Value = 1.28
$ws.Range("P10").Value = 2.02
$ws.Range("Q10").Value = 1.87
$ws.Range("R10").Value = 1.42
$ws.Range("S10").Value = 3.2
$ws.Range("T10").Value = 2.06
$ws.Range("U10").Value = 1.76
$ws.Range("V10").Value = 1.08
$ws.Range("W10").Value = 3.3
$ws.Range("X10").Value = 17
$ws.Range("AA10").Value = 420
$ws.Range("AB10").Value = 8
$ws.Range("AD10").Value = 40
$ws.Range("AE10").Value = 190
$ws.Range("AG10").Value = 10
$ws.Range("AJ10").Value = 12
$ws.Range("AO10").Value = 250
$ws.Range("F11").Value = 3.9
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 3.25
$ws.Range("L11").Value = 1.32
$ws.Range("N11").Value = 4
$ws.Range("Q11").Value = 1.55
$ws.Range("S11").Value = 2.38
$ws.Range("V11").Value = 2
$ws.Range("Y11").Value = 980
$ws.Range("AD11").Value = 13
$ws.Range("F12").Value = 2.72
$ws.Range("G12").Value = 2.92
$ws.Range("H12").Value = 2.84
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3.15
$ws.Range("K12").Value = 3.3
$ws.Range("L12").Value = 1.56
$ws.Range("M12").Value = 1.11
$ws.Range("N12").Value = 2.84
$ws.Range("O12").Value = 1.49
$ws.Range("P12").Value = 1.6
$ws.Range("Q12").Value = 2.48
$ws.Range("R12").Value = 1.21
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 2.02
$ws.Range("U12").Value = 1.81
$ws.Range("V12").Value = 1.47
$ws.Range("W12").Value = 1.52
$ws.Range("X12").Value = 10.5
$ws.Range("Y12").Value = 9.800000000000001
$ws.Range("Z12").Value = 19.5
$ws.Range("AB12").Value = 9.800000000000001
$ws.Range("AC12").Value = 7.2
$ws.Range("AF12").Value = 16.5
$ws.Range("AH12").Value = 22
$ws.Range("AJ12").Value = 48
$ws.Range("AK12").Value = 42
$ws.Range("AL12").Value = 65
$ws.Range("AM12").Value = 170
$ws.Range("AN12").Value = 46
$ws.Range("AO12").Value = 1000
